# 1. Slide 5's table: switch its table style from the default
#    "Table_0" style ({EA1ED201-9E51-446C-937C-2DA6A24097D2}) to
#    {01CC63C0-A8B5-4FDA-B09B-9262B1BD95F8}.
$p = $ppt.ActivePresentation
$tableSlide = $p.Slides.Item(5)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{01CC63C0-A8B5-4FDA-B09B-9262B1BD95F8}")
    }
}

# 2. Re-colour the deck's theme so the master's theme carries the
#    "Office" palette (the deck was originally themed with the
#    "Integral" / Red Violet palette; the edit swaps it back to the
#    plain Office colours) while the table-style change above stays
#    local to slide 5.
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Item(1).RGB  = 0          # dk1      -> 000000
$tcs.Item(2).RGB  = 16777215   # lt1      -> FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      -> 44546A
$tcs.Item(4).RGB  = 15132391   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  -> ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  -> FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  -> 4472C4
$tcs.Item(10).RGB = 4697456    # accent6  -> 70AD47
$tcs.Item(11).RGB = 12673797   # hlink    -> 0563C1
$tcs.Item(12).RGB = 7491477    # folHlink -> 954F72
